$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.910.26'
$ws.Range('E2').Value = '  +0.97%  '
$ws.Range('D3').Value = '2.358.88'
$ws.Range('E3').Value = '  +4.52%  '
$ws.Range('E4').Value = '  -0.22%  '
$ws.Range('B5').Value = 'XRP'
$ws.Range('C5').Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.660'
$ws.Range('E5').Value = '  +3.14%  '
$ws.Range('B6').Value = 'BNB'
$ws.Range('C6').Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '235.16'
$ws.Range('E6').Value = '  +1.58%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '73.45'
$ws.Range('E7').Value = '  +14.12%  '
$ws.Range('E8').Value = '  -0.11%  '
$ws.Range('E9').Value = '  +21.88%  '
$ws.Range('E10').Value = '  +3.41%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '28.46'
$ws.Range('E11').Value = '  +8.12%  '
$ws.Range('B12').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C12').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D12').Value = '2.711.19'
$ws.Range('E12').Value = '  +4.58%  '
$ws.Range('B13').Value = 'TRON'
$ws.Range('C13').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.106'
$ws.Range('E13').Value = '  +2.44%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '16.83'
$ws.Range('E14').Value = '  +12.48%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.63'
$ws.Range('E15').Value = '  +9.74%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.885'
$ws.Range('E16').Value = '  +7.60%  '
$ws.Range('D17').Value = '2.353.23'
$ws.Range('E17').Value = '  +4.27%  '
$ws.Range('D18').Value = '43.764.14'
$ws.Range('E18').Value = '  +0.88%  '
$ws.Range('E19').Value = '  +4.59%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '76.05'
$ws.Range('E20').Value = '  +4.33%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.32'
$ws.Range('E21').Value = '  +3.95%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '251.02'
$ws.Range('E22').Value = '  +1.84%  '
$ws.Range('B23').Value = 'Dai'
$ws.Range('C23').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.00'
$ws.Range('E23').Value = '  +0.06%  '
$ws.Range('B24').Value = 'WEMIXToken'
$ws.Range('C24').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '3.79'
$ws.Range('E24').Value = '  -2.95%  '
$ws.Range('E25').Value = '  +1.36%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '10.28'
$ws.Range('E26').Value = '  +5.93%  '
$ws.Range('E27').Value = '  -1.13%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '22.51'
$ws.Range('E28').Value = '  +4.13%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '172.86'
$ws.Range('E29').Value = '  -0.53%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.55'
$ws.Range('E30').Value = '  +8.52%  '
$ws.Range('E31').Value = '  +1.67%  '
$ws.Range('E32').Value = '  +4.88%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.16'
$ws.Range('E33').Value = '  +4.54%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0707'
$ws.Range('E34').Value = '  +4.34%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.12'
$ws.Range('E35').Value = '  +4.58%  '
$ws.Range('E36').Value = '  +4.12%  '
$ws.Range('E37').Value = '  +8.03%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '6.43'
$ws.Range('E38').Value = '  +0.36%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0264'
$ws.Range('E39').Value = '  +5.64%  '
$ws.Range('E40').Value = '  +13.44%  '
$ws.Range('E41').Value = '  -0.05%  '
$ws.Range('E42').Value = '  +1.38%  '
$ws.Range('E43').Value = '  +10.06%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '99.03'
$ws.Range('E44').Value = '  +2.65%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.21'
$ws.Range('E45').Value = '  +3.15%  '
$ws.Range('E46').Value = '  +2.37%  '
$ws.Range('E47').Value = '  -1.27%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.180'
$ws.Range('E48').Value = '  +13.16%  '
$ws.Range('D49').Value = '1.438.82'
$ws.Range('E49').Value = '  +0.98%  '
$ws.Range('B50').Value = 'TerraClassic'
$ws.Range('C50').Value = 'https://coinranking.com/coin/AaQUAs2Mc+terraclassic-lunc'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.000205'
$ws.Range('E50').Value = '  -0.30%  '
$ws.Range('B51').Value = 'NEARProtocol'
$ws.Range('C51').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.30'
$ws.Range('E51').Value = '  +1.79%  '
